$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.950.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.553.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("E6").Value = "  +0.35%  "

$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.14"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.06%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0858"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.20%  "

$ws.Range("E12").Value = "  +0.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.555.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.50%  "

$ws.Range("E15").Value = "  +1.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.956.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0698"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.61%  "

$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.45%  "

$ws.Range("E24").Value = "  +0.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.87%  "

$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.85%  "

$ws.Range("E28").Value = "  +1.12%  "

$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0467"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.41%  "

$ws.Range("E32").Value = "  +0.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.425.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.53%  "

$ws.Range("E35").Value = "  +3.38%  "

$ws.Range("E36").Value = "  +1.46%  "

$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("E38").Value = "  +0.71%  "

$ws.Range("E39").Value = "  +0.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.809"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.57%  "

$ws.Range("E41").Value = "  +3.63%  "

$ws.Range("E42").Value = "  -0.17%  "

$ws.Range("E43").Value = "  +4.70%  "

$ws.Range("E44").Value = "  +0.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.689.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0520"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0954"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.74%  "
